$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update training data values from -1 to 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("A3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("C5").Value = 0

# Update the selected cell on the sheet view
$ws.Range("D5").Select()
